# Updated cryptos list values (price + 1h volume) to match refreshed data feed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.060.10'
$ws.Range("E2").Value = '  +3.05%  '
$ws.Range("D3").Value = '1.654.89'
$ws.Range("E3").Value = '  +3.77%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.68%  '
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +1.73%  '
$ws.Range("E9").Value = '  +1.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0864'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.35%  '
$ws.Range("D12").Value = '1.886.85'
$ws.Range("E12").Value = '  +3.72%  '
$ws.Range("D13").Value = '1.648.04'
$ws.Range("E14").Value = '  +1.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.520'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.88%  '
$ws.Range("D17").Value = '27.029.67'
$ws.Range("E17").Value = '  +3.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '237.89'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.48%  '
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("E22").Value = '  +4.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.26%  '
$ws.Range("E28").Value = '  +1.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.34%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("E31").Value = '  +1.38%  '
$ws.Range("D32").Value = '1.526.51'
$ws.Range("E32").Value = '  +4.28%  '
$ws.Range("E33").Value = '  +2.81%  '
$ws.Range("E34").Value = '  +3.28%  '
$ws.Range("E35").Value = '  +8.22%  '
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.887'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.07%  '
$ws.Range("E39").Value = '  +2.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.70%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.49%  '
$ws.Range("E43").Value = '  +3.76%  '
$ws.Range("D44").Value = '1.794.93'
$ws.Range("E44").Value = '  +3.61%  '
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.920'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.12'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("E49").Value = '  +3.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0504'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("E51").Value = '  +3.10%  '
